# Added medial I and Is under vowels
# Splits the combined "N" and "Y" medial rows on the "Medial" sheet into
# separate coda / medial rows, moving the "A, E, I, O, U, V, AE, IE, IY"
# contextual-suffix entries that used to live at the bottom of the sheet
# (rows 76 & 77) up next to their matching "N"/"Y" medial rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Medial")

# Remove the two trailing rows that held the stand-alone coda entries
# (delete from the bottom up so row numbers above stay valid).
$ws.Rows.Item(77).Delete()
$ws.Rows.Item(76).Delete()

# --- N: split row 24 ("N" / "N.med") into N.coda (24) + N.med (new row 25)
$ws.Rows.Item(25).Insert()
$ws.Range("C24").Value = "N.coda"
$ws.Range("B25").Value = "N"
$ws.Range("C25").Value = "N.med"
$ws.Range("D25").Value = "A, E, I, O, U, V, AE, IE, IY"

# --- Y: split row 51 ("Y" / "Y.med") into Y.coda (51) + Y.med (new row 52)
$ws.Rows.Item(52).Insert()
$ws.Range("C51").Value = "Y.coda"
$ws.Range("B52").Value = "Y"
$ws.Range("C52").Value = "Y.med"
$ws.Range("D52").Value = "A, E, I, O, U, V, AE, IE, IY"

# Make "Medial" the active sheet with D52 selected (it was previously the
# "Final" sheet that had focus).
$ws.Activate()
$ws.Range("D52").Select()
